$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.281.85"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "1.662.88"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.73%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.56"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5314"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2637"
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06365"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.51"
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07851"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.554"
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").Value = "1.671.29"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "1.892.35"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "0.0₅8177"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.66"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.010"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.658"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.52"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.056"
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.10"
$ws.Range("E24").Value = "  +3.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1220"
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.231"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.486"
$ws.Range("E28").Value = "  +3.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05884"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.585"
$ws.Range("E31").Value = "  +2.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.304"
$ws.Range("E32").Value = "  +2.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.613"
$ws.Range("E33").Value = "  +4.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9584"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.818"
$ws.Range("E35").Value = "  +2.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.424"
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5812"
$ws.Range("E37").Value = "  +3.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01614"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.896"
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8545"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("D42").Value = "1.048.47"
$ws.Range("E42").Value = "  +3.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.18"
$ws.Range("E43").Value = "  +2.84%  "
$ws.Range("D44").Value = "1.805.07"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.27"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").Value = "0.0₈107"
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.014"
$ws.Range("E47").Value = "  +1.06%  "
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.976"
$ws.Range("E49").Value = "  +3.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05163"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.442"
